$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed (semana 4-7)
$ws.Range("B5").Value = 634
$ws.Range("B6").Value = 434
$ws.Range("B7").Value = 458
$ws.Range("B8").Value = 495

# Add new rows for semana 9-12 (2025 week 11 data update)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 499

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 582

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 556

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 3
